$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F2"
$ws.Range("C2").Value = "F2rl2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3845463333333334
$ws.Range("H2").Value = 1.153639
$ws.Range("I2").Value = 0.1984850200147207
$ws.Range("J2").Value = 0.1984850200147207
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.163436
$ws.Range("N2").Value = 0.490308
$ws.Range("O2").Value = 0.1789685622010547
$ws.Range("P2").Value = 0.1789685622010547
$ws.Range("Q2").Value = 0.06284871453466667
$ws.Range("R2").Value = 0.565638430812
$ws.Range("S2").Value = 0.03552257865048214
$ws.Range("T2").Value = 0.03552257865048213

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F2"
$ws.Range("C3").Value = "F2rl2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3845463333333334
$ws.Range("H3").Value = 1.153639
$ws.Range("I3").Value = 0.1984850200147207
$ws.Range("J3").Value = 0.1984850200147207
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7497746666666667
$ws.Range("N3").Value = 2.249324
$ws.Range("O3").Value = 0.8210314377989453
$ws.Range("P3").Value = 0.8210314377989452
$ws.Range("Q3").Value = 0.2883230988928889
$ws.Range("R3").Value = 2.594907890036
$ws.Range("S3").Value = 0.1629624413642386
$ws.Range("T3").Value = 0.1629624413642385

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F2"
$ws.Range("C4").Value = "F2rl2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.110828
$ws.Range("H4").Value = 3.332484
$ws.Range("I4").Value = 0.5733580031870772
$ws.Range("J4").Value = 0.5733580031870772
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.163436
$ws.Range("N4").Value = 0.490308
$ws.Range("O4").Value = 0.1789685622010547
$ws.Range("P4").Value = 0.1789685622010547
$ws.Range("Q4").Value = 0.181549285008
$ws.Range("R4").Value = 1.633943565072
$ws.Range("S4").Value = 0.102613057456859
$ws.Range("T4").Value = 0.1026130574568589

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "F2"
$ws.Range("C5").Value = "F2rl2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.110828
$ws.Range("H5").Value = 3.332484
$ws.Range("I5").Value = 0.5733580031870772
$ws.Range("J5").Value = 0.5733580031870772
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7497746666666667
$ws.Range("N5").Value = 2.249324
$ws.Range("O5").Value = 0.8210314377989453
$ws.Range("P5").Value = 0.8210314377989452
$ws.Range("Q5").Value = 0.832870693424
$ws.Range("R5").Value = 7.495836240816001
$ws.Range("S5").Value = 0.4707449457302182
$ws.Range("T5").Value = 0.4707449457302182

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "F2"
$ws.Range("C6").Value = "F2rl2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.442033
$ws.Range("H6").Value = 1.326099
$ws.Range("I6").Value = 0.2281569767982021
$ws.Range("J6").Value = 0.2281569767982021
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.163436
$ws.Range("N6").Value = 0.490308
$ws.Range("O6").Value = 0.1789685622010547
$ws.Range("P6").Value = 0.1789685622010547
$ws.Range("Q6").Value = 0.072244105388
$ws.Range("R6").Value = 0.6501969484919999
$ws.Range("S6").Value = 0.04083292609371364
$ws.Range("T6").Value = 0.04083292609371363

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "F2"
$ws.Range("C7").Value = "F2rl2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.442033
$ws.Range("H7").Value = 1.326099
$ws.Range("I7").Value = 0.2281569767982021
$ws.Range("J7").Value = 0.2281569767982021
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7497746666666667
$ws.Range("N7").Value = 2.249324
$ws.Range("O7").Value = 0.8210314377989453
$ws.Range("P7").Value = 0.8210314377989452
$ws.Range("Q7").Value = 0.3314251452306666
$ws.Range("R7").Value = 2.982826307076
$ws.Range("S7").Value = 0.1873240507044885
$ws.Range("T7").Value = 0.1873240507044884

# Remove now-unused rows 8-10 (MuSCs sending cluster rows removed in new data)
$ws.Rows("8:10").Delete()